$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 1 de Junio de 2020 a las 19:05"

# Row 4: Estados Unidos
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 1843304
$ws.Range("C4").Value = 6134
$ws.Range("D4").Value = 601049
$ws.Range("E4").Value = 1135922
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 138
$ws.Range("H4").Value = 106333

# Row 9: Italia
$ws.Range("A9").Value = "Italia"
$ws.Range("B9").Value = 233197
$ws.Range("C9").Value = 200
$ws.Range("D9").Value = 158355
$ws.Range("E9").Value = 41367
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 60
$ws.Range("H9").Value = 33475

# Row 10: India
$ws.Range("A10").Value = "India"
$ws.Range("B10").Value = 197808
$ws.Range("C10").Value = 7199
$ws.Range("D10").Value = 95415
$ws.Range("E10").Value = 96790
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 195
$ws.Range("H10").Value = 5603

# Row 12: Alemania
$ws.Range("A12").Value = "Alemania"
$ws.Range("B12").Value = 183606
$ws.Range("C12").Value = 112
$ws.Range("D12").Value = 165900
$ws.Range("E12").Value = 9096
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 5
$ws.Range("H12").Value = 8610

# Row 13: Turquia
$ws.Range("A13").Value = "Turquia"
$ws.Range("B13").Value = 164769
$ws.Range("C13").Value = 827
$ws.Range("D13").Value = 128947
$ws.Range("E13").Value = 31259
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 23
$ws.Range("H13").Value = 4563

# Row 14: Peru
$ws.Range("A14").Value = "Peru"
$ws.Range("B14").Value = 164476
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = 67208
$ws.Range("E14").Value = 92762
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 4506

# Row 37: Irlanda
$ws.Range("A37").Value = "Irlanda"
$ws.Range("B37").Value = 25062
$ws.Range("C37").Value = 72
$ws.Range("D37").Value = 22089
$ws.Range("E37").Value = 1321
$ws.Range("F37").Value = 0
$ws.Range("G37").Value = 0
$ws.Range("H37").Value = 1652

# Row 58: Argelia
$ws.Range("A58").Value = "Argelia"
$ws.Range("B58").Value = 9513
$ws.Range("C58").Value = 119
$ws.Range("D58").Value = 5894
$ws.Range("E58").Value = 2958
$ws.Range("F58").Value = 0
$ws.Range("G58").Value = 8
$ws.Range("H58").Value = 661

# Row 59: Armenia
$ws.Range("A59").Value = "Armenia"
$ws.Range("B59").Value = 9492
$ws.Range("C59").Value = 210
$ws.Range("D59").Value = 3402
$ws.Range("E59").Value = 5951
$ws.Range("F59").Value = 0
$ws.Range("G59").Value = 8
$ws.Range("H59").Value = 139

# Row 128: Republica del Chad
$ws.Range("A128").Value = "Republica del Chad"
$ws.Range("B128").Value = 790
$ws.Range("C128").Value = 12
$ws.Range("D128").Value = 539
$ws.Range("E128").Value = 185
$ws.Range("F128").Value = 0
$ws.Range("G128").Value = 1
$ws.Range("H128").Value = 66

# Row 142: Cabo Verde
$ws.Range("A142").Value = "Cabo Verde"
$ws.Range("B142").Value = 458
$ws.Range("C142").Value = 23
$ws.Range("D142").Value = 193
$ws.Range("E142").Value = 261
$ws.Range("F142").Value = 0
$ws.Range("G142").Value = 0
$ws.Range("H142").Value = 4

# Row 143: Uganda
$ws.Range("A143").Value = "Uganda"
$ws.Range("B143").Value = 457
$ws.Range("C143").Value = 40
$ws.Range("D143").Value = 72
$ws.Range("E143").Value = 385
$ws.Range("F143").Value = 0
$ws.Range("G143").Value = 0
$ws.Range("H143").Value = 0

# Row 144: Estado de Palestina
$ws.Range("A144").Value = "Estado de Palestina"
$ws.Range("B144").Value = 449
$ws.Range("C144").Value = 1
$ws.Range("D144").Value = 372
$ws.Range("E144").Value = 74
$ws.Range("F144").Value = 0
$ws.Range("G144").Value = 0
$ws.Range("H144").Value = 3

# Row 145: Taiwan
$ws.Range("A145").Value = "Taiwan"
$ws.Range("B145").Value = 443
$ws.Range("C145").Value = 1
$ws.Range("D145").Value = 427
$ws.Range("E145").Value = 9
$ws.Range("F145").Value = 0
$ws.Range("G145").Value = 0
$ws.Range("H145").Value = 7

# Row 146: Togo
$ws.Range("A146").Value = "Togo"
$ws.Range("B146").Value = 442
$ws.Range("C146").Value = 0
$ws.Range("D146").Value = 211
$ws.Range("E146").Value = 218
$ws.Range("F146").Value = 0
$ws.Range("G146").Value = 0
$ws.Range("H146").Value = 13
